# CSCE155N Final Project Written Submission - edit script
# Applies the text and layout changes described by the target diff.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                                      $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
    }
}

# --- Paragraph 1 (introduction) ---------------------------------------

Replace-Text "a process that is used to solve" "a process that’s used to solve"

Replace-Text "industry. Chemical Engineers work" "industry. Chemical engineers work"

Replace-Text "both everyday life such as water treatment and more specialized such as nuclear power plants." "both everyday life and specialized applications."

Replace-Text "While optimizing these processes we are often" "While optimizing these processes, we are often"

$old1 = "in a substance, this may be atoms, molecules or ions. In one mole there are 6.022 x 10"
$old1 += "23"
$old1 += " particles. Converting"
Replace-Text $old1 "in a substance - atoms, molecules, or ions. Converting"

Replace-Text "This gives the grams (or another unit of measurement) of each reactant or product." "This gives the mass of each reactant or product."

Replace-Text "This program will both speed up the trivial step and will eliminate much of the human error that can occur." "This program speeds up the trivial step and eliminates much of the human error that can occur."

# --- "calculations as chemical engineering majors." ---------------------

Replace-Text "calculations as chemical engineering majors. " "calculations. "

# --- Quinn's words paragraph --------------------------------------------

Replace-Text "a very ambitious project.” However, as a result, our understanding of GUI and its applications in MATLAB increased greatly. Additionally, we now have more knowledge of how to manipulate GUI" "a very ambitious project.” As a result, our understanding of GUI and its applications in MATLAB increased greatly and we have more knowledge of how to manipulate GUI"

# --- GUI programming paragraph -------------------------------------------

Replace-Text "Many issues in the Chemical industry can be solved through multiple simple steps and calculations. While it is necessary" "Many issues in the chemical industry are be solved through multiple simple steps and calculations. While it’s necessary"

Replace-Text "many of the steps could and should be completed by a computer program to reduce human error and increase workplace efficiency. GUI programs are also very user friendly due to the display options that programmers have. Some specific examples" "many could and should be completed by a computer program to reduce human error and increase efficiency. Some specific examples"

# --- Closing paragraph -----------------------------------------------------

Replace-Text "e is able to convert mass and mole fractions using an easy-to-understand interface" "e converts mole fractions to mass fractions using an easy-to-understand interface"

# Remove the trailing empty paragraph (right before the section break).
$paraCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($paraCount)
if ($lastPara.Range.Text.Trim().Length -eq 0) {
    $lastPara.Range.Delete()
}

# --- Header: merge "Quinn " + "Lanik" runs, drop proofErr spell markers ---

Replace-Text "Topic Proposal Approved by Quinn Lanik" "Topic Proposal Approved by Quinn Lanik"

# --- Section page margins --------------------------------------------------

$d.PageSetup.TopMargin = 72
$d.PageSetup.BottomMargin = 72
$d.PageSetup.LeftMargin = 54
$d.PageSetup.RightMargin = 54
